$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05881766666666666
$ws.Range("H2").Value = 0.176453
$ws.Range("M2").Value = 1.090291
$ws.Range("N2").Value = 3.270873
$ws.Range("O2").Value = 0.02878663098006733
$ws.Range("P2").Value = 0.02878663098006733
$ws.Range("Q2").Value = 0.06412837260766666
$ws.Range("R2").Value = 0.5771553534689999
$ws.Range("S2").Value = 0.02878663098006733
$ws.Range("T2").Value = 0.02878663098006733

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05881766666666666
$ws.Range("H3").Value = 0.176453
$ws.Range("O3").Value = 0.7239320554917256
$ws.Range("P3").Value = 0.7239320554917257
$ws.Range("Q3").Value = 1.61271336786
$ws.Range("R3").Value = 14.51442031074
$ws.Range("S3").Value = 0.7239320554917256
$ws.Range("T3").Value = 0.7239320554917257

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05881766666666666
$ws.Range("H4").Value = 0.176453
$ws.Range("O4").Value = 0.2472813135282071
$ws.Range("P4").Value = 0.2472813135282071
$ws.Range("Q4").Value = 0.550871973307
$ws.Range("R4").Value = 4.957847759763
$ws.Range("S4").Value = 0.2472813135282071
$ws.Range("T4").Value = 0.2472813135282071
